$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...respecting the speaker and also other people..."
#    -> "...respecting the speaker and other people..."
#    (this also absorbs/removes the gramStart/gramEnd proofErr pair that
#    wrapped "and also")
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("respecting the speaker and also other", $true, $false, $false, $false, $false, $true, 1, $false, "respecting the speaker and other", 2) | Out-Null

# re-split so "and" lives in its own run again, matching the original layout
$r1b = $d.Content
$r1b.Find.Execute("and", $true, $true) | Out-Null
$r1b.Font.Bold = 1
$r1b.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) Gameplay paragraph:
#    "...causing another students to talk inside the auditorium. So the
#    player needs..."
#    -> "...causing another student to talk inside the auditorium. So, the
#    player needs..."
#    (removes the gramStart/gramEnd proofErr pairs around "another" and "So")
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("causing another students to talk", $true, $false, $false, $false, $false, $true, 1, $false, "causing another student to talk", 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute("auditorium. So the player", $true, $false, $false, $false, $false, $true, 1, $false, "auditorium. So, the player", 2) | Out-Null

# re-split "student" into its own run
$r4 = $d.Content
$r4.Find.Execute("student", $true, $true) | Out-Null
$r4.Font.Bold = 1
$r4.Font.Bold = 0

# re-split "So," into its own run
$r5 = $d.Content
$r5.Find.Execute("So,", $true, $false) | Out-Null
$r5.Font.Bold = 1
$r5.Font.Bold = 0

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the document to wrap the
#    Gameplay paragraph just edited above (this mirrors Word recording the
#    location of the most recent edit).
# ---------------------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
if ($oldBm -ne $null) {
    $oldBm.Delete()
}

$targetPara = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("The game play is simple")) {
        $targetPara = $i
        break
    }
}
if ($targetPara -ne -1) {
    $pRange = $d.Paragraphs($targetPara).Range
    $bmRange = $d.Range($pRange.Start, $pRange.End - 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 4) "The player could tap the noisy students, an animation of a students "
#    split off the second "students" into its own run (proofErr spell-check
#    wrap around it in the source edit).
# ---------------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.Execute("an animation of a students", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitStart = $r6.End - 8
$splitEnd = $r6.End
$r6b = $d.Range($splitStart, $splitEnd)
$r6b.Font.Bold = 1
$r6b.Font.Bold = 0
